# Remove redundant hypercapnia columns
# - Rename header D1/E1 (drop "other_hypercap_threshold", shift "pco2_threshold_any" into D1,
#   add new "unknown_hypercap_threshold" into E1)
# - Refresh pivoted data: insert a new "Administrative" category row, drop the
#   "Symptom – Musculoskeletal" category row, and update all aggregated values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("D1").Value = "pco2_threshold_any"
$ws.Range("E1").Value = "unknown_hypercap_threshold"

# --- Data rows (rows 2-12) ---
$data = @(
    @("Administrative",              2.28, 1.51, 1.94,  2.27, 1.54),
    @("Diseases (patient-stated)",   6.83, 2.98, 5.73,  5.56, 4.45),
    @("Injuries & adverse effects", 21.07, 7.97, 16.49, 16.41, 11.95),
    @("Other",                       5.58, 5.75, 6.05,  7.07, 5.85),
    @("Symptom – Circulatory",       6.46, 6.71, 7.57,  9.09, 8.300000000000001),
    @("Symptom – Digestive",         9.74, 7.21, 9.880000000000001, 10.1, 9.210000000000001),
    @("Symptom – General",           2.76, 4.14, 3.39,  3.79, 3.68),
    @("Symptom – Nervous",          10.45, 11.25, 11.83, 11.62, 13.12),
    @("Symptom – Respiratory",      28.67, 49.87, 32.61, 29.55, 39.35),
    @("Symptom – Skin/Hair/Nails",   2.22, 1.36, 1.71,  2.27, 0.9399999999999999),
    @("Uncodable/Unknown",           3.93, 1.26, 2.79,  2.27, 1.63)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r++
}
